$d = $word.ActiveDocument

# Helper: replace the visible text of a Range (a Cell.Range or Paragraph.Range)
# without disturbing the trailing cell-mark / paragraph-mark, and without
# relying on Find.Execute (whose Replace scope is not honoured by this
# runtime - it always operates on the whole document).
function Set-RangeText($rng, $newText) {
    $target = $d.Range($rng.Start, $rng.End - 1)
    $target.Text = $newText
}

# Update the date line at the top of the document.
Set-RangeText $d.Paragraphs.Item(1).Range "2025-05-29 Thursday"

$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    Set-RangeText $t.Cell($row, $col).Range $newText
}

# Row 1 (table row index 1)
Set-CellText 1 1 "238×8=1904"
Set-CellText 1 2 "234×4=936"
Set-CellText 1 3 "979×9=8811"
Set-CellText 1 4 "456×5=2280"
Set-CellText 1 5 "182×3=546"

# Row 5 (table row index 5)
Set-CellText 5 1 "129×6=774"
Set-CellText 5 2 "343×6=2058"
Set-CellText 5 3 "356×2=712"
Set-CellText 5 4 "280×3=840"
Set-CellText 5 5 "935×5=4675"

# Row 10 (table row index 10)
Set-CellText 10 1 "174×3=522"
Set-CellText 10 2 "345×6=2070"
Set-CellText 10 3 "144×8=1152"
Set-CellText 10 4 "167×2=334"
Set-CellText 10 5 "234×7=1638"

# Row 15 (table row index 15)
Set-CellText 15 1 "120×2=240"
Set-CellText 15 2 "233×8=1864"
Set-CellText 15 3 "737×3=2211"
Set-CellText 15 4 "292×8=2336"
Set-CellText 15 5 "803×9=7227"

# Row 20 (table row index 20)
Set-CellText 20 1 "262×2=524"
Set-CellText 20 2 "990×4=3960"
Set-CellText 20 3 "523×5=2615"
Set-CellText 20 4 "214×3=642"
Set-CellText 20 5 "785×4=3140"

Write-Output "Done"
